$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Step 1: "Complex" heading paragraph gets a new paragraph added right
# after it: a lead-in sentence plus the _GoBack bookmark that used to sit
# at the end of the "Average" paragraph above. Do this one first since
# "Complex" is the very last paragraph in the body - editing it can't shift
# the position of anything earlier in the document (namely the "Average"
# paragraph we edit in step 2). ---

$complexPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Style.NameLocal -eq "Heading 2" -and $cand.Range.Text.Trim() -eq "Complex") {
        $complexPara = $cand
        break
    }
}

# Remove the old (mis-placed) _GoBack bookmark so the new one added below
# becomes id 0 cleanly, with no stale duplicate left behind.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# InsertXML on a collapsed range at the tail of a paragraph replaces that
# paragraph instead of appending after it, so re-emit the heading
# paragraph's own content alongside the new paragraph and replace the
# whole original paragraph range with both.
$newComplexXml = (
    "<w:p $wns><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:t>Complex</w:t></w:r></w:p>" +
    "<w:p $wns>" +
        "<w:r><w:t xml:space='preserve'>The complex methodology illustrates how to build a </w:t></w:r>" +
        "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
        "<w:bookmarkEnd w:id='0'/>" +
    "</w:p>"
)
[void]$complexPara.Range.InsertXML($newComplexXml)

# --- Step 2: The "Average" methodology paragraph (the one that used to end
# with the bookmark) gains a new trailing sentence, then a blank paragraph,
# then a new paragraph describing population of the fact table (with a
# proofed "methodology," run). Re-locate it fresh now that step 1 has
# already run. ---

$avgPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("The average methodology")) {
        $avgPara = $cand
        break
    }
}

$newAvgXml = (
    "<w:p $wns>" +
        "<w:r><w:t>The average methodology is similar to the simple methodology in that the names of the</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> generated dimensional tables are prefixed with " + [char]8220 + "DIM_" + [char]8221 + " and the resulting fact table is prefixed with the letters " + [char]8220 + "FACT_" + [char]8221 + ".  Beyond that, the methodology is more complex in that it illustrates how the user of the application can generate custom dimensional tables using a powerful table builder screen in the application (Screenshot #2).</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>  Similar to the simple methodology, a single fact table is generated from the custom dimensional tables.</w:t></w:r>" +
    "</w:p>" +
    "<w:p/>" +
    "<w:p>" +
        "<w:r><w:t xml:space='preserve'>Population of the dimensional tables is automatic, but just as in the simple </w:t></w:r>" +
        "<w:proofErr w:type='gramStart'/>" +
        "<w:r><w:t>methodology,</w:t></w:r>" +
        "<w:proofErr w:type='gramEnd'/>" +
        "<w:r><w:t xml:space='preserve'> population of the fact table requires manual script writing and execution.</w:t></w:r>" +
    "</w:p>"
)
[void]$avgPara.Range.InsertXML($newAvgXml)
